# Refresh the cryptocurrency price/volume table (GitHub Actions data pull).
# Each entry below is a cell address together with its new text value, taken
# from the authoritative diff of the workbook. Every value in columns D and E
# is stored as text in the sheet (prices use "." as a thousands separator in
# some rows, e.g. "98.285.31", and percentages keep their padding spaces,
# e.g. "  -0.55%  "), so it must round-trip as text rather than be
# reinterpreted as a number. Values that "look like" a plain number (e.g.
# "9.00", "0.150", "0.0000268") would otherwise be auto-converted by Excel
# and lose trailing zeros or be rewritten in scientific notation, so those
# are entered with a leading apostrophe to force text, and the quote-prefix
# style this introduces is immediately reset back to "Normal" so the cell's
# appearance/formatting is unaffected.

$updates = @(
    ,@("D2", "98.285.31")
    ,@("E2", "  -0.55%  ")
    ,@("D3", "3.417.54")
    ,@("E3", "  +2.03%  ")
    ,@("E4", "  +0.02%  ")
    ,@("D5", "255.48")
    ,@("E5", "  -1.71%  ")
    ,@("D6", "665.78")
    ,@("E6", "  +2.54%  ")
    ,@("E7", "  -5.60%  ")
    ,@("D8", "0.437")
    ,@("E8", "  -6.38%  ")
    ,@("D9", "1.06")
    ,@("E9", "  -2.23%  ")
    ,@("E10", "  +0.00%  ")
    ,@("D11", "3.415.52")
    ,@("E11", "  +2.05%  ")
    ,@("E12", "  +3.14%  ")
    ,@("D13", "42.52")
    ,@("E13", "  -2.85%  ")
    ,@("D14", "6.53")
    ,@("E14", "  +16.71%  ")
    ,@("D15", "97.908.06")
    ,@("E15", "  -2.05%  ")
    ,@("D16", "0.0000268")
    ,@("E16", "  -0.71%  ")
    ,@("D17", "4.060.48")
    ,@("E17", "  +1.91%  ")
    ,@("D18", "8.93")
    ,@("E18", "  +18.90%  ")
    ,@("D19", "3.427.65")
    ,@("E19", "  +2.15%  ")
    ,@("D20", "0.558")
    ,@("E20", "  +28.64%  ")
    ,@("D21", "17.69")
    ,@("E21", "  +4.55%  ")
    ,@("D22", "11.25")
    ,@("E22", "  +8.90%  ")
    ,@("D23", "3.47")
    ,@("E23", "  -3.56%  ")
    ,@("D24", "511.97")
    ,@("E24", "  -5.07%  ")
    ,@("D25", "0.0000207")
    ,@("E25", "  -2.98%  ")
    ,@("D26", "6.74")
    ,@("E26", "  +7.85%  ")
    ,@("D27", "101.23")
    ,@("E27", "  -1.62%  ")
    ,@("D28", "12.95")
    ,@("E28", "  +1.62%  ")
    ,@("D29", "3.594.89")
    ,@("E29", "  +1.90%  ")
    ,@("D30", "0.150")
    ,@("E30", "  +0.60%  ")
    ,@("D31", "11.66")
    ,@("E31", "  +5.74%  ")
    ,@("D32", "0.198")
    ,@("E32", "  +2.72%  ")
    ,@("E33", "  -0.09%  ")
    ,@("D34", "2.41")
    ,@("E34", "  +15.55%  ")
    ,@("D35", "0.999")
    ,@("E35", "  -0.06%  ")
    ,@("E36", "  +7.42%  ")
    ,@("E37", "  +1.82%  ")
    ,@("D38", "1.54")
    ,@("E38", "  +16.61%  ")
    ,@("D39", "8.01")
    ,@("E39", "  +2.46%  ")
    ,@("D40", "540.06")
    ,@("E40", "  +4.04%  ")
    ,@("E41", "  -1.15%  ")
    ,@("E42", "  +0.07%  ")
    ,@("D43", "0.878")
    ,@("E43", "  +6.22%  ")
    ,@("D44", "24.72")
    ,@("E44", "  +0.01%  ")
    ,@("B45", "Cosmos")
    ,@("C45", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom")
    ,@("D45", "9.00")
    ,@("E45", "  +14.80%  ")
    ,@("B46", "ImmutableX")
    ,@("C46", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx")
    ,@("D46", "1.76")
    ,@("E46", "  +17.70%  ")
    ,@("B47", "MantraDAO")
    ,@("C47", "https://coinranking.com/coin/cTdD8lD-6+mantradao-om")
    ,@("D47", "3.78")
    ,@("E47", "  -0.10%  ")
    ,@("B48", "Filecoin")
    ,@("C48", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil")
    ,@("D48", "5.84")
    ,@("E48", "  +14.63%  ")
    ,@("B49", "VeChain")
    ,@("C49", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet")
    ,@("D49", "0.0432")
    ,@("E49", "  +0.56%  ")
    ,@("D50", "3.29")
    ,@("E50", "  -2.96%  ")
    ,@("D51", "53.93")
    ,@("E51", "  +8.97%  ")
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($update in $updates) {
    $address = $update[0]
    $value = $update[1]
    $cell = $ws.Range($address)
    if ($value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        # Looks like a plain number: force text so trailing zeros / exponent
        # notation survive the round-trip, the way typing an apostrophe before
        # a number does in the Excel UI, then drop the quote-prefix style it
        # introduces so the cell's formatting is left exactly as it was.
        $cell.Value = "'" + $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
